# Apply scheduled runner updates to Titan_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 286
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 286
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 286
$ws.Range("M9").Value = $null
$ws.Range("N9").Value = -624
# Row 21
$ws.Range("H21").Value = 30499.285
$ws.Range("J21").Value = 28999.23
$ws.Range("L21").Value = 28999.23
$ws.Range("N21").Value = -29935.23
# Row 23
$ws.Range("H23").Value = 30499.285
$ws.Range("J23").Value = 28999.23
$ws.Range("L23").Value = 28999.23
$ws.Range("N23").Value = -29467.23
# Row 29
$ws.Range("H29").Value = 5507.857
$ws.Range("J29").Value = 5422.2
$ws.Range("L29").Value = 16266.6
$ws.Range("N29").Value = -16828.6
# Row 38
$ws.Range("H38").Value = 2709.6
$ws.Range("I38").Value = 565.6667
$ws.Range("J38").Value = 3628.4285
$ws.Range("K38").Value = 1697.0001
$ws.Range("L38").Value = 10885.2855
$ws.Range("M38").Value = -1325.0001
$ws.Range("N38").Value = -11629.2855
# Row 58
$ws.Range("H58").Value = 1700.8334
$ws.Range("I58").Value = 235
$ws.Range("J58").Value = 3166.6667
$ws.Range("K58").Value = 705
$ws.Range("L58").Value = 9500.000100000001
$ws.Range("M58").Value = -555
$ws.Range("N58").Value = -9800.000100000001
# Row 87
$ws.Range("H87").Value = 33966.668
$ws.Range("J87").Value = 33966.668
$ws.Range("L87").Value = 33966.668
$ws.Range("N87").Value = -36462.668
# Row 90
$ws.Range("H90").Value = 33966.668
$ws.Range("J90").Value = 33966.668
$ws.Range("L90").Value = 101900.004
$ws.Range("N90").Value = -114380.004
# Row 98
$ws.Range("H98").Value = 256134.34
$ws.Range("I98").Value = 274665.25
$ws.Range("J98").Value = 2878.6667
$ws.Range("K98").Value = 274665.25
$ws.Range("L98").Value = 2878.6667
$ws.Range("M98").Value = -273167.25
$ws.Range("N98").Value = -5874.6667
# Row 122
$ws.Range("H122").Value = 256134.34
$ws.Range("I122").Value = 274665.25
$ws.Range("J122").Value = 2878.6667
$ws.Range("K122").Value = 823995.75
$ws.Range("L122").Value = 8636.000100000001
$ws.Range("M122").Value = -821545.75
$ws.Range("N122").Value = -13536.0001
# Row 137
$ws.Range("H137").Value = 21740330
$ws.Range("I137").Value = 27778610
$ws.Range("J137").Value = 2523
$ws.Range("K137").Value = 83335830
$ws.Range("L137").Value = 7569
$ws.Range("M137").Value = -83333280
$ws.Range("N137").Value = -12669

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1607
$ws.Range("I45").Value = 1237.25
$ws.Range("K45").Value = 1237.25
$ws.Range("M45").Value = -860.25
# Row 61
$ws.Range("H61").Value = 3264.5293
$ws.Range("I61").Value = 1992.0588
$ws.Range("J61").Value = 4537
$ws.Range("K61").Value = 1992.0588
$ws.Range("L61").Value = 4537
$ws.Range("M61").Value = -1780.0588
$ws.Range("N61").Value = -4961
# Row 74
$ws.Range("H74").Value = 5805.9644
$ws.Range("I74").Value = 1465.2
$ws.Range("J74").Value = 10814.538
$ws.Range("K74").Value = 1465.2
$ws.Range("L74").Value = 10814.538
$ws.Range("M74").Value = -591.2
$ws.Range("N74").Value = -12562.538
# Row 77
$ws.Range("H77").Value = 5805.9644
$ws.Range("I77").Value = 1465.2
$ws.Range("J77").Value = 10814.538
$ws.Range("K77").Value = 7326
$ws.Range("L77").Value = 54072.69
$ws.Range("M77").Value = -2958
$ws.Range("N77").Value = -62808.69
# Row 123
$ws.Range("H123").Value = 28619.334
$ws.Range("J123").Value = 28619.334
$ws.Range("L123").Value = 28619.334
$ws.Range("N123").Value = -38419.334
# Row 132
$ws.Range("H132").Value = 3407.5386
$ws.Range("I132").Value = 3009.3076
$ws.Range("J132").Value = 3805.7693
$ws.Range("K132").Value = 9027.9228
$ws.Range("L132").Value = 11417.3079
$ws.Range("M132").Value = -6497.9228
$ws.Range("N132").Value = -16477.3079
# Row 136
$ws.Range("H136").Value = 3264.5293
$ws.Range("I136").Value = 1992.0588
$ws.Range("J136").Value = 4537
$ws.Range("K136").Value = 5976.1764
$ws.Range("L136").Value = 13611
$ws.Range("M136").Value = -3426.1764
$ws.Range("N136").Value = -18711

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2479.975
$ws.Range("I134").Value = 1974.4839
$ws.Range("J134").Value = 4221.1113
$ws.Range("K134").Value = 5923.4517
$ws.Range("L134").Value = 12663.3339
$ws.Range("M134").Value = -3388.4517
$ws.Range("N134").Value = -17733.3339

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1236.3125
$ws.Range("I31").Value = 1030.2858
$ws.Range("K31").Value = 1030.2858
$ws.Range("M31").Value = -735.2858000000001
# Row 34
$ws.Range("H34").Value = 1236.3125
$ws.Range("I34").Value = 1030.2858
$ws.Range("K34").Value = 1030.2858
$ws.Range("M34").Value = -828.2858000000001
# Row 58
$ws.Range("H58").Value = 1725.8125
$ws.Range("I58").Value = 1172
$ws.Range("J58").Value = 4125.6665
$ws.Range("K58").Value = 1172
$ws.Range("L58").Value = 4125.6665
$ws.Range("M58").Value = -969
$ws.Range("N58").Value = -4531.6665
# Row 132
$ws.Range("H132").Value = 2597.7407
$ws.Range("I132").Value = 2203.5715
$ws.Range("J132").Value = 3977.3333
$ws.Range("K132").Value = 6610.7145
$ws.Range("L132").Value = 11931.9999
$ws.Range("M132").Value = -4080.7145
$ws.Range("N132").Value = -16991.9999
# Row 134
$ws.Range("H134").Value = 1935.8667
$ws.Range("J134").Value = 4820.4546
$ws.Range("L134").Value = 14461.3638
$ws.Range("N134").Value = -19531.3638
# Row 136
$ws.Range("H136").Value = 1725.8125
$ws.Range("I136").Value = 1172
$ws.Range("J136").Value = 4125.6665
$ws.Range("K136").Value = 3516
$ws.Range("L136").Value = 12376.9995
$ws.Range("M136").Value = -966
$ws.Range("N136").Value = -17476.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 13514313
$ws.Range("I113").Value = 589.5
$ws.Range("J113").Value = 15152340
$ws.Range("K113").Value = 1768.5
$ws.Range("L113").Value = 45457020
$ws.Range("M113").Value = 401.5
$ws.Range("N113").Value = -45461360

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2811.238
$ws.Range("I132").Value = 2639.7917
$ws.Range("J132").Value = 3039.8333
$ws.Range("K132").Value = 7919.375100000001
$ws.Range("L132").Value = 9119.499899999999
$ws.Range("M132").Value = -5389.375100000001
$ws.Range("N132").Value = -14179.4999

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550
# Row 136
$ws.Range("H136").Value = 3509.848
$ws.Range("I136").Value = 2081.8484
$ws.Range("J136").Value = 7134.769
$ws.Range("K136").Value = 6245.5452
$ws.Range("L136").Value = 21404.307
$ws.Range("M136").Value = -3695.5452
$ws.Range("N136").Value = -26504.307

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 65293.562
$ws.Range("I122").Value = 113266.336
$ws.Range("J122").Value = 3614.2856
$ws.Range("K122").Value = 339799.008
$ws.Range("L122").Value = 10842.8568
$ws.Range("M122").Value = -337349.008
$ws.Range("N122").Value = -15742.8568
# Row 123
$ws.Range("H123").Value = 25459.215
$ws.Range("J123").Value = 25459.215
$ws.Range("L123").Value = 25459.215
$ws.Range("N123").Value = -35259.215
# Row 132
$ws.Range("H132").Value = 11630231
$ws.Range("J132").Value = 1051
$ws.Range("L132").Value = 3153
$ws.Range("N132").Value = -8213

